# 自动更新Excel文件 - 2026-02-08 23:22:25
# Advance the tracker's "today" reference by one day (2026-02-08 -> 2026-02-09):
#  - For each data row (2..99), column E holds remaining days, column F holds
#    the start date (yyyymmdd integer), column D holds the total day count.
#  - If remaining (E) was 1, the cycle has completed: reset F to the new
#    "today" (20260209) and reset E back to the full total (D).
#  - Otherwise simply decrement E by 1 (F stays put).
#  - Rows whose F value isn't a parseable yyyymmdd date (e.g. row 36) are
#    left untouched, mirroring the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newToday = 20260209

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    # NOTE: read via .Value2 (bare .Value reads are unreliable in this
    # COM shim); writes use .Value2 as well for symmetry.
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    # Validate F looks like a real yyyymmdd date (8 digits, valid calendar date).
    $fStr = [string]([int]$fVal)
    $isValidDate = $false
    if ($fStr.Length -eq 8) {
        $y = [int]$fStr.Substring(0, 4)
        $m = [int]$fStr.Substring(4, 2)
        $d = [int]$fStr.Substring(6, 2)
        if ($m -ge 1 -and $m -le 12 -and $d -ge 1 -and $d -le 31) {
            $isValidDate = $true
        }
    }

    if (-not $isValidDate) {
        continue
    }

    if ($eVal -eq 1) {
        $eCell.Value2 = $dVal
        $fCell.Value2 = $newToday
    } else {
        $eCell.Value2 = $eVal - 1
    }
}
